# Insert a new weekly record row at row 83 (Feria Lagunitas de Puerto Montt - Acelga),
# pushing the existing rows 83-123 down to 84-124.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(83).Insert()

$ws.Cells.Item(83, 1).Value = 4
$ws.Cells.Item(83, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(83, 3).Value = "Los Lagos"
$ws.Cells.Item(83, 4).Value = 44523
$ws.Cells.Item(83, 5).Value = 10
$ws.Cells.Item(83, 6).Value = 100112009
$ws.Cells.Item(83, 7).Value = "Acelga"
$ws.Cells.Item(83, 8).Value = "Sin especificar"
$ws.Cells.Item(83, 9).Value = "Primera"
$ws.Cells.Item(83, 10).Value = 200
$ws.Cells.Item(83, 11).Value = 3500
$ws.Cells.Item(83, 12).Value = 3500
$ws.Cells.Item(83, 13).Value = 3500
$ws.Cells.Item(83, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(83, 15).Value = "Región del Maule"
$ws.Cells.Item(83, 16).Value = 875
$ws.Cells.Item(83, 17).Value = 4
$ws.Cells.Item(83, 18).Value = "Hortaliza"

$ws.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
